# Events.xlsx - "Added hooks in for choosing path to retreat in counterattack"
#
# Insert a new row 110 on the "Events" sheet holding a new event "e099b"
# (a retreat-choice prompt used when a counterattack retreat has more than
# one valid destination hex). All rows that were previously 110-120 shift
# down to 111-121.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 110 (and everything below it) down by one row.
$ws.Rows("110").Insert()

# Column A: short event id used elsewhere in the workbook as a lookup key.
$ws.Range("A110").Value2 = "e099b"

# Column B: the rich "event text" markup shown to the player.
$body = "<Bold>e099b Retreat Choice in Counterattack Retreat</Bold> " + "`n" + `
    "<InlineUIContainer><Button Content='r11.33' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  " + "`n" + `
    "<InlineUIContainer><Button Content='r20.45' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>     " + "`n" + `
    "<LineBreak/><LineBreak/>" + "`n" + `
    "Must retreat to an adjacent hex toward the start area. Since two or more choices exist, choose one of the highlighted areas on the movement board to continue." + "`n" + `
    "<LineBreak/><LineBreak/>" + "`n" + `
    "                        <InlineUIContainer><Image Name='Sherman1' Height='200' Width='325'></Image></InlineUIContainer>  "

$ws.Range("B110").Value2 = $body

# Match the row height used by the row this new one was modelled on.
$ws.Rows("110").RowHeight = 105

# Keep the active selection pinned where the author's cursor ended up.
$ws.Range("B109").Select()
